$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16512555771186125"
$wb.Worksheets.Item(2).Name = "NB_TO-1651255579562095"
$wb.Worksheets.Item(3).Name = "RS_TO-1651255579562095"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512555796083114"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16512555796801126"

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16512555770856018.csv"
$ws1.Range("B3").Value = "GNG_stims-1651255577099602.csv"
$ws1.Range("B4").Value = "go_stims-1651255577101604.csv"
$ws1.Range("B5").Value = "GNG_stims-16512555771156027.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16512555790847518.csv"
$ws2.Range("B3").Value = "ZB-match_1-16512555773371646.csv"
$ws2.Range("B4").Value = "OB-16512555783602843.csv"
$ws2.Range("B5").Value = "ZB-match_3-16512555771871672.csv"
$ws2.Range("B6").Value = "OB-16512555779491065.csv"
$ws2.Range("B7").Value = "OB-16512555781141076.csv"
$ws2.Range("B8").Value = "TB-16512555786291716.csv"
$ws2.Range("B9").Value = "TB-16512555795464668.csv"
$ws2.Range("B10").Value = "ZB-match_0-16512555776561093.csv"

# --- Sheet 3: RS (no data changes, name already updated above) ---

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16512555795781183.csv"
$ws4.Range("B3").Value = "ZM_stims-1651255579562095.csv"
$ws4.Range("B4").Value = "MM_stims-16512555795940323.csv"
$ws4.Range("B5").Value = "ZM_stims-16512555795791183.csv"
$ws4.Range("B6").Value = "MM_stims-16512555796083114.csv"
$ws4.Range("B7").Value = "ZM_stims-16512555795951068.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16512555796640954.csv"
$ws5.Range("B3").Value = "SAT_stims-1651255579636343.csv"
$ws5.Range("B4").Value = "vSAT_stims-16512555796524415.csv"
$ws5.Range("B5").Value = "SAT_stims-16512555796128814.csv"
